$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 7-17 with revised AntalTest (B) / AntalOmikron (C) values.
# Ratio (D) already holds a shared formula "=100*C/B" that will recalculate automatically.

$updates = @(
    @{ Row = 7;  B = 3919; C = 3 },
    @{ Row = 8;  B = 3835; C = 11 },
    @{ Row = 9;  B = 4803; C = 12 },
    @{ Row = 10; B = 5181; C = 24 },
    @{ Row = 11; B = 4267; C = 77 },
    @{ Row = 12; B = 4294; C = 62 },
    @{ Row = 13; B = 4946; C = 75 },
    @{ Row = 14; B = 5089; C = 111 },
    @{ Row = 15; B = 4995; C = 167 },
    @{ Row = 16; B = 6762; C = 337 },
    @{ Row = 17; B = 7039; C = 530 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 2).Value = $u.B
    $ws.Cells.Item($u.Row, 3).Value = $u.C
}

# Extend the shared Ratio formula down through the new rows (18:20) and add the
# new dates / counts for 2021-12-08, 2021-12-09 and 2021-12-10.

$newRows = @(
    @{ Row = 18; Date = "'2021-12-08"; B = 6204; C = 643 },
    @{ Row = 19; Date = "'2021-12-09"; B = 4223; C = 501 },
    @{ Row = 20; Date = "'2021-12-10"; B = 2912; C = 381 }
)

foreach ($n in $newRows) {
    $ws.Cells.Item($n.Row, 1).Value = $n.Date
    $ws.Cells.Item($n.Row, 2).Value = $n.B
    $ws.Cells.Item($n.Row, 3).Value = $n.C
    $ws.Cells.Item($n.Row, 4).Formula = "=100*C$($n.Row)/B$($n.Row)"
}

$ws.Range("C28").Select()

$wb.Save()
